# Add 2022-Q4 data:
#  - Insert a new worksheet "2022-Q4" right before "2022-Q3" (becomes the
#    2nd tab, right after "总计"), populated with the quarterly fund-holder
#    breakdown.
#  - Insert a new leading row in "总计" summarizing the 2022-Q4 totals,
#    pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for 2022-Q4, shifting the rest down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Rows(2).Insert()

$idxSummaryCell = $summary.Cells.Item(2, 1)
$idxSummaryCell.Value = 0
$idxSummaryCell.Font.Bold = $true
$idxSummaryCell.HorizontalAlignment = -4108
$idxSummaryCell.VerticalAlignment = -4160
$idxSummaryCell.Borders.LineStyle = 1
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 18
$summary.Cells.Item(2, 4).Value = 4.81

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet, inserted before the existing "2022-Q3" tab.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q4.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @(0,  "519702", "交银趋势优先混合A",             "99.56", "81.36", "2.03", "2.0211", 10),
    @(1,  "014038", "交银启诚混合A",                 "24.58", "81.04", "2.82", "0.6932", 4),
    @(2,  "013430", "交银趋势优先混合C",             "20.19", "81.36", "2.03", "0.4099", 10),
    @(3,  "001487", "宝盈优势产业灵活配置混合A",     "10.62", "94.48", "3.24", "0.3441", 6),
    @(4,  "001128", "宝盈新兴产业灵活配置混合A",     "9.13",  "94.44", "3.73", "0.3405", 6),
    @(5,  "001877", "宝盈国家安全沪港深股票A",       "6.58",  "94.32", "3.95", "0.2599", 8),
    @(6,  "013895", "宝盈成长精选混合A",             "8.59",  "94.68", "3.00", "0.2577", 8),
    @(7,  "014039", "交银启诚混合C",                 "7.22",  "81.04", "2.82", "0.2036", 4),
    @(8,  "012771", "宝盈优势产业灵活配置混合C",     "3.19",  "94.48", "3.24", "0.1034", 6),
    @(9,  "013896", "宝盈成长精选混合C",             "3.06",  "94.68", "3.00", "0.0918", 8),
    @(10, "012815", "宝盈新兴产业灵活配置混合C",     "1.06",  "94.44", "3.73", "0.0395", 6),
    @(11, "002378", "建信弘利灵活配置混合A",         "0.82",  "92.97", "3.27", "0.0268", 10),
    @(12, "006072", "民生加银创新成长混合A",         "0.40",  "91.73", "3.33", "0.0133", 4),
    @(13, "013613", "宝盈国家安全沪港深股票C",       "0.13",  "94.32", "3.95", "0.0051", 8),
    @(14, "003855", "汇安丰华灵活配置混合C",         "0.18",  "45.67", "2.22", "0.0040", 9),
    @(15, "014929", "民生加银创新成长混合C",         "0.01",  "91.73", "3.33", "0.0003", 4),
    @(16, "017194", "建信弘利灵活配置混合C",         "0.00",  "92.97", "3.27", $null,    10),
    @(17, "003854", "汇安丰华灵活配置混合A",         "0.00",  "45.67", "2.22", $null,    9)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $q4.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $codeCell = $q4.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[1]

    $q4.Cells.Item($r, 3).Value = $row[2]

    for ($c = 4; $c -le 6; $c++) {
        $textCell = $q4.Cells.Item($r, $c)
        $textCell.NumberFormat = "@"
        $textCell.Value = $row[$c - 1]
    }

    $gValue = $row[6]
    if ($null -eq $gValue) {
        $q4.Cells.Item($r, 7).Value = 0
    }
    else {
        $gCell = $q4.Cells.Item($r, 7)
        $gCell.NumberFormat = "@"
        $gCell.Value = $gValue
    }

    $q4.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

Write-Output "2022-Q4 sheet and summary row added"
